$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237, shifting existing rows 237-244 down to 238-245.
$ws.Rows.Item(237).Insert()

# The new row 237 receives the data that used to be in row 236 (before the D/J update below).
$ws.Cells.Item(237, 1).Value = 10
$ws.Cells.Item(237, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(237, 3).Value = "La Araucanía"
$ws.Cells.Item(237, 4).Value = 44167
$ws.Cells.Item(237, 5).Value = 9
$ws.Cells.Item(237, 6).Value = 100114014
$ws.Cells.Item(237, 7).Value = "Betarraga"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 35
$ws.Cells.Item(237, 11).Value = 8000
$ws.Cells.Item(237, 12).Value = 8000
$ws.Cells.Item(237, 13).Value = 8000
$ws.Cells.Item(237, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(237, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(237, 16).Value = 667
$ws.Cells.Item(237, 17).Value = 12
$ws.Cells.Item(237, 18).Value = "Hortaliza"

# Row 236 is updated with a new weekly record (new date and volume).
$ws.Cells.Item(236, 4).Value = 44448
$ws.Cells.Item(236, 10).Value = 125

Write-Host "Row inserted and values updated"
